# Apply calibration updates to the "Frac of Start Year Cap Ret per Unit Net Loss" workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "FoSYCRpUNL": update retirement-fraction values (column B) ---
$ws = $wb.Worksheets.Item("FoSYCRpUNL")

# Most plant types move from 0.5% (or 0.1%, in biomass's case) up to 1.25%.
# Petroleum (B12) and natural gas peaker (B13) are intentionally left untouched
# because they are called out as exceptions (see About sheet note added below).
$rowsToRaise = 2,4,5,6,7,8,9,10,11,14,15,16,17,18,19,20,21,22,23,24
foreach ($r in $rowsToRaise) {
    $ws.Cells.Item($r, 2).Value = 0.0125
}

# natural gas steam turbine keeps a very low fraction, but it is recalibrated too
$ws.Cells.Item(3, 2).Value = 0.0001

$ws.Range("B3").Select() | Out-Null

# --- Sheet "About": add explanatory note about low-retirement plant types ---
$about = $wb.Worksheets.Item("About")
$about.Cells.Item(8, 1).Value = "Certain plant types are less prone to economic retirement because they are maintained for local reliability purposes."
$about.Cells.Item(10, 1).Value = "These includes: natural gas steam turbines, natural gas peakers, and petroleum plants. For these plant types we set the "
$about.Cells.Item(11, 1).Value = "retirement fraction very low."

$about.Range("A12").Select() | Out-Null
